$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to be interpreted/stored as text so that numeric-looking
    # strings (prices, percentages) keep their exact formatting (leading/
    # trailing zeros, sign, "%") instead of being coerced into a Double.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Restore the default "Normal" style so we do not leave a stray custom
    # number-format style behind (the source cells have no explicit style).
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "301.00"
Set-TextValue $ws.Range("E2") "-3.04%"
Set-TextValue $ws.Range("D3") "35.46"
Set-TextValue $ws.Range("E3") "-0.24%"
Set-TextValue $ws.Range("D4") "5.057"
Set-TextValue $ws.Range("E4") "-0.94%"
Set-TextValue $ws.Range("D5") "0.07983"
Set-TextValue $ws.Range("E5") "-2.60%"
Set-TextValue $ws.Range("D6") "1.895"
Set-TextValue $ws.Range("E6") "-7.86%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.049"
Set-TextValue $ws.Range("E7") "-1.78%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "7.749"
Set-TextValue $ws.Range("E8") "-2.49%"
Set-TextValue $ws.Range("D9") "0.9275"
Set-TextValue $ws.Range("E9") "0.20%"
Set-TextValue $ws.Range("D10") "0.1443"
Set-TextValue $ws.Range("D11") "0.1899"
Set-TextValue $ws.Range("E11") "-0.68%"
Set-TextValue $ws.Range("D12") "0.09029"
Set-TextValue $ws.Range("E12") "-2.47%"
Set-TextValue $ws.Range("D13") "0.03474"
Set-TextValue $ws.Range("E13") "-4.83%"
Set-TextValue $ws.Range("D14") "0.09851"
Set-TextValue $ws.Range("E14") "-0.57%"
Set-TextValue $ws.Range("D15") "0.001397"
Set-TextValue $ws.Range("E15") "-2.37%"
Set-TextValue $ws.Range("D16") "0.005815"
Set-TextValue $ws.Range("E16") "-0.57%"
Set-TextValue $ws.Range("D17") "3.535"
Set-TextValue $ws.Range("E17") "1.80%"
Set-TextValue $ws.Range("D18") "2.992"
Set-TextValue $ws.Range("E18") "1.16%"
Set-TextValue $ws.Range("D19") "0.3426"
Set-TextValue $ws.Range("E19") "0.89%"
Set-TextValue $ws.Range("D20") "0.1302"
Set-TextValue $ws.Range("E20") "-0.44%"
Set-TextValue $ws.Range("D21") "5.036"
Set-TextValue $ws.Range("E21") "-1.16%"
Set-TextValue $ws.Range("E22") "8.47%"
Set-TextValue $ws.Range("D23") "0.04486"
Set-TextValue $ws.Range("E23") "-0.95%"
Set-TextValue $ws.Range("D24") "0.001213"
Set-TextValue $ws.Range("E24") "-1.10%"
Set-TextValue $ws.Range("D25") "0.004766"
Set-TextValue $ws.Range("E25") "-0.89%"
Set-TextValue $ws.Range("E26") "-1.62%"
Set-TextValue $ws.Range("D27") "0.0003023"
Set-TextValue $ws.Range("E27") "-31.94%"
Set-TextValue $ws.Range("D39") "0.01831"
Set-TextValue $ws.Range("E39") "-6.98%"
Set-TextValue $ws.Range("D40") "0.04753"
Set-TextValue $ws.Range("E40") "-2.83%"
Set-TextValue $ws.Range("D41") "0.01051"
Set-TextValue $ws.Range("E41") "15.98%"
Set-TextValue $ws.Range("D42") "0.007329"
Set-TextValue $ws.Range("E42") "-3.69%"
Set-TextValue $ws.Range("D43") "0.1327"
Set-TextValue $ws.Range("E43") "-4.11%"
Set-TextValue $ws.Range("E44") "-3.57%"
Set-TextValue $ws.Range("D45") "0.01090"
Set-TextValue $ws.Range("E45") "-6.34%"
Set-TextValue $ws.Range("D46") "0.00006223"
Set-TextValue $ws.Range("E46") "-5.00%"
Set-TextValue $ws.Range("E47") "0.09%"
Set-TextValue $ws.Range("D48") "64.67"
Set-TextValue $ws.Range("E48") "-64.06%"
Set-TextValue $ws.Range("E49") "10.77%"
Set-TextValue $ws.Range("D50") "0.00002099"
Set-TextValue $ws.Range("E50") "0.09%"
Set-TextValue $ws.Range("D51") "0.0001999"
Set-TextValue $ws.Range("E51") "0.09%"
